$wb = $excel.ActiveWorkbook

# Sheet 1 (展览) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 6).Value = 446
$ws1.Cells.Item(5, 6).Value = 1341
$ws1.Cells.Item(6, 6).Value = 7677
$ws1.Cells.Item(9, 6).Value = 2098
$ws1.Cells.Item(10, 6).Value = 8465
$ws1.Cells.Item(13, 6).Value = 66
$ws1.Cells.Item(14, 6).Value = 5670
$ws1.Cells.Item(16, 6).Value = 2630
$ws1.Cells.Item(17, 6).Value = 1144
$ws1.Cells.Item(18, 6).Value = 4596
$ws1.Cells.Item(20, 6).Value = 405
$ws1.Cells.Item(23, 6).Value = 538
$ws1.Cells.Item(24, 6).Value = 3563
$ws1.Cells.Item(25, 6).Value = 60
$ws1.Cells.Item(27, 6).Value = 26
$ws1.Cells.Item(28, 6).Value = 4
$ws1.Cells.Item(29, 6).Value = 3058
$ws1.Cells.Item(30, 6).Value = 45
$ws1.Cells.Item(31, 6).Value = 115
$ws1.Cells.Item(32, 6).Value = 350
$ws1.Cells.Item(33, 6).Value = 129
$ws1.Cells.Item(34, 6).Value = 317
$ws1.Cells.Item(35, 6).Value = 609
$ws1.Cells.Item(36, 6).Value = 663
$ws1.Cells.Item(38, 6).Value = 884
$ws1.Cells.Item(39, 6).Value = 1850
$ws1.Cells.Item(40, 6).Value = 48
$ws1.Cells.Item(42, 6).Value = 21
$ws1.Cells.Item(43, 6).Value = 2966
$ws1.Cells.Item(44, 6).Value = 6
$ws1.Cells.Item(45, 6).Value = 2291
$ws1.Cells.Item(46, 6).Value = 11
$ws1.Cells.Item(47, 6).Value = 31

# Sheet 2 (演出) - column F updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 129
$ws2.Cells.Item(9, 6).Value = 122

# Sheet 4 (全部类型) - column F updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(5, 6).Value = 1341
$ws4.Cells.Item(6, 6).Value = 7677
$ws4.Cells.Item(9, 6).Value = 2098
$ws4.Cells.Item(10, 6).Value = 8465
$ws4.Cells.Item(12, 6).Value = 66
$ws4.Cells.Item(13, 6).Value = 5670
$ws4.Cells.Item(15, 6).Value = 2630
$ws4.Cells.Item(16, 6).Value = 1144
$ws4.Cells.Item(17, 6).Value = 4596
$ws4.Cells.Item(18, 6).Value = 405
$ws4.Cells.Item(22, 6).Value = 129
$ws4.Cells.Item(23, 6).Value = 538
$ws4.Cells.Item(25, 6).Value = 3563
$ws4.Cells.Item(26, 6).Value = 60
$ws4.Cells.Item(28, 6).Value = 26
$ws4.Cells.Item(29, 6).Value = 4
$ws4.Cells.Item(30, 6).Value = 3058
$ws4.Cells.Item(31, 6).Value = 350
$ws4.Cells.Item(32, 6).Value = 129
$ws4.Cells.Item(33, 6).Value = 317
$ws4.Cells.Item(35, 6).Value = 610
$ws4.Cells.Item(36, 6).Value = 663
$ws4.Cells.Item(38, 6).Value = 884
$ws4.Cells.Item(40, 6).Value = 1850
$ws4.Cells.Item(41, 6).Value = 48
$ws4.Cells.Item(43, 6).Value = 21
$ws4.Cells.Item(44, 6).Value = 2966
$ws4.Cells.Item(45, 6).Value = 2291
$ws4.Cells.Item(46, 6).Value = 31
$ws4.Cells.Item(48, 6).Value = 122
